# Generated by Katalon AI
# Adds two new data rows (6 & 7) to the "AI Generated" sheet and widens
# columns D, E, F slightly to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen columns D (4), E (5), F (6) ---------------------------------
# NOTE: this runtime's ColumnWidth setter stores width + 5/6 characters
# (it bakes in a fixed padding offset), so back the requested width off
# by 5/6 here to land on the exact target character widths (18/17/16).
$ws.Columns.Item(4).ColumnWidth = 18 - 5/6
$ws.Columns.Item(5).ColumnWidth = 17 - 5/6
$ws.Columns.Item(6).ColumnWidth = 16 - 5/6

# --- plain text values for row 6 & 7 (columns A-F) ---------------------
$ws.Cells.Item(6, 1).Value = "f03 address"
$ws.Cells.Item(6, 2).Value = "03 city"
$ws.Cells.Item(6, 3).Value = " "
$ws.Cells.Item(6, 4).Value = "f03 first"
$ws.Cells.Item(6, 5).Value = "f03 last"
$ws.Cells.Item(6, 6).Value = " "

$ws.Cells.Item(7, 1).Value = "15 address"
$ws.Cells.Item(7, 2).Value = "15 HCM city"
$ws.Cells.Item(7, 3).Value = " "
$ws.Cells.Item(7, 4).Value = "15 silicon first"
$ws.Cells.Item(7, 5).Value = "15 silicon last"
$ws.Cells.Item(7, 6).Value = "{{ip_address}}"

# --- numeric-looking values for G/H must stay TEXT (keep leading zeros) -
# Writing "0303" directly via .Value gets auto-converted to the number 303
# (dropping the leading zero), since it looks like a number. Stage the
# values in scratch cells using a leading quote (forces text), copy them,
# and paste-special "values" into the real destination cells - that keeps
# the target cells as genuine text without carrying the scratch
# formatting along, matching the original inline-string cells.
$scratch = $ws.Range($ws.Cells.Item(500, 20), $ws.Cells.Item(501, 21))
$ws.Cells.Item(500, 20).Value = "'0303"
$ws.Cells.Item(500, 21).Value = "'030303"
$ws.Cells.Item(501, 20).Value = "'1515"
$ws.Cells.Item(501, 21).Value = "'151515"

$scratch.Copy()
$ws.Range("G6:H7").PasteSpecial(-4163)
$scratch.Clear()
